# Daily attendance processing
# For every row in the "Recorded By" column (G), whenever the cell holds
# multiple comma-separated recorder names/emails, reverse the order of
# the entries (single-entry cells are left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text

    if ($text -eq "") {
        continue
    }

    $parts = $text -split ", "
    $count = $parts.Count

    if ($count -gt 1) {
        $reversed = @()
        for ($i = $count - 1; $i -ge 0; $i--) {
            $reversed += $parts[$i]
        }
        $cell.Value = [string]::Join(", ", $reversed)
    }
}
